$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, applied as text to preserve original string formatting
# (prices/volume values are stored as text in this sheet, e.g. "0.0000140", "2.495.23").
$updates = @{
    'D2' = '59.387.82'
    'E2' = '  -4.14%  '
    'D3' = '2.495.23'
    'E3' = '  -3.20%  '
    'D4' = '1.00'
    'E4' = '  +0.05%  '
    'D5' = '536.53'
    'E5' = '  -2.42%  '
    'D6' = '144.38'
    'E6' = '  -6.90%  '
    'E7' = '  -0.36%  '
    'D8' = '0.575'
    'E8' = '  -3.46%  '
    'D9' = '2.530.94'
    'E9' = '  -2.01%  '
    'D10' = '0.101'
    'E10' = '  -3.39%  '
    'E11' = '  -2.66%  '
    'D12' = '5.58'
    'E12' = '  +0.12%  '
    'D13' = '0.354'
    'E13' = '  -3.46%  '
    'D14' = '2.941.40'
    'E14' = '  -3.06%  '
    'D15' = '24.12'
    'E15' = '  -6.02%  '
    'D16' = '59.270.95'
    'E16' = '  -4.21%  '
    'D17' = '0.0000140'
    'E17' = '  -3.63%  '
    'D18' = '2.515.06'
    'E18' = '  -2.60%  '
    'D19' = '11.36'
    'E19' = '  -2.29%  '
    'D20' = '4.31'
    'E20' = '  -5.58%  '
    'D21' = '325.04'
    'E21' = '  -3.97%  '
    'D22' = '0.996'
    'E22' = '  -0.31%  '
    'D23' = '5.77'
    'E23' = '  -4.60%  '
    'D24' = '61.29'
    'E24' = '  -3.54%  '
    'D25' = '0.443'
    'E25' = '  -10.44%  '
    'D26' = '0.163'
    'E26' = '  -2.93%  '
    'D27' = '2.614.87'
    'E27' = '  -3.06%  '
    'D28' = '0.987'
    'E28' = '  -1.20%  '
    'D29' = '7.83'
    'E29' = '  -4.39%  '
    'D30' = '6.98'
    'E30' = '  -4.98%  '
    'D31' = '0.0₃0782'
    'E31' = '  -6.81%  '
    'D32' = '1.26'
    'E32' = '  -6.17%  '
    'D33' = '1.79'
    'E33' = '  -5.64%  '
    'E34' = '  -0.28%  '
    'D35' = '158.08'
    'E35' = '  -2.90%  '
    'D36' = '1.45'
    'E36' = '  +1.11%  '
    'D37' = '18.58'
    'E37' = '  -3.43%  '
    'D38' = '4.46'
    'E38' = '  -8.50%  '
    'D39' = '1.62'
    'E39' = '  -9.94%  '
    'D40' = '5.91'
    'E40' = '  -2.72%  '
    'D41' = '310.29'
    'E41' = '  -5.77%  '
    'D42' = '36.85'
    'E42' = '  -1.98%  '
    'D43' = '3.70'
    'E43' = '  -6.29%  '
    'D44' = '0.825'
    'E44' = '  -9.73%  '
    'D45' = '0.995'
    'E45' = '  -0.19%  '
    'D46' = '0.599'
    'E46' = '  -1.37%  '
    'D47' = '10.79'
    'E47' = '  -1.34%  '
    'D48' = '125.02'
    'D49' = '0.0933'
    'E49' = '  -3.56%  '
    'D50' = '18.73'
    'E50' = '  -4.25%  '
    'D51' = '0.0520'
    'E51' = '  -5.18%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text number-format so Excel does not auto-convert numeric-looking
    # strings (e.g. "1.00", "0.0000140") into actual numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Drop the temporary text format so the cell keeps its original (default) style.
    $cell.ClearFormats()
}
